$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H28").Value = 2484.9
$ws.Range("I28").Value = 284.33334
$ws.Range("J28").Value = 3428
$ws.Range("K28").Value = 284.33334
$ws.Range("L28").Value = 3428
$ws.Range("M28").Value = 200.66666
$ws.Range("N28").Value = -4398

$ws.Range("H87").Value = 49900
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 49900
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 49900
$ws.Range("N87").Value = -52396

$ws.Range("H90").Value = 49900
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 49900
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 149700
$ws.Range("N90").Value = -162180

$ws.Range("H92").Value = 370
$ws.Range("I92").Value = 370
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 370
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 878
$ws.Range("N92").ClearContents()

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()

$ws.Range("H95").Value = 0
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H96").Value = 1550
$ws.Range("I96").Value = 400
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1200
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = 173

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H132").Value = 1134.5834
$ws.Range("I132").Value = 1134.5834
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3403.7502
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -873.7501999999999

$ws.Range("H137").Value = 1833.9375
$ws.Range("I137").Value = 1073.1666
$ws.Range("J137").Value = 2290.4
$ws.Range("K137").Value = 3219.4998
$ws.Range("L137").Value = 6871.200000000001
$ws.Range("M137").Value = -669.4998000000001
$ws.Range("N137").Value = -11971.2

$ws.Range("H138").Value = 3429.4443
$ws.Range("I138").Value = 3336.6191
$ws.Range("J138").Value = 3559.4
$ws.Range("K138").Value = 10009.8573
$ws.Range("L138").Value = 10678.2
$ws.Range("M138").Value = -4869.8573
$ws.Range("N138").Value = -20958.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2630.8044
$ws.Range("I32").Value = 2090.8555
$ws.Range("J32").Value = 7610.3335
$ws.Range("K32").Value = 2090.8555
$ws.Range("L32").Value = 7610.3335
$ws.Range("M32").Value = -1803.8555
$ws.Range("N32").Value = -8184.3335

$ws.Range("H61").Value = 3730.182
$ws.Range("I61").Value = 1237.25
$ws.Range("J61").Value = 5154.7144
$ws.Range("K61").Value = 1237.25
$ws.Range("L61").Value = 5154.7144
$ws.Range("M61").Value = -1025.25
$ws.Range("N61").Value = -5578.7144

$ws.Range("H74").Value = 877.06665
$ws.Range("I74").Value = 834.1818
$ws.Range("J74").Value = 995
$ws.Range("K74").Value = 834.1818
$ws.Range("L74").Value = 995
$ws.Range("M74").Value = 39.81820000000005
$ws.Range("N74").Value = -2743

$ws.Range("H77").Value = 877.06665
$ws.Range("I77").Value = 834.1818
$ws.Range("J77").Value = 995
$ws.Range("K77").Value = 4170.909
$ws.Range("L77").Value = 4975
$ws.Range("M77").Value = 197.0910000000003
$ws.Range("N77").Value = -13711

$ws.Range("H122").Value = 1795.9231
$ws.Range("I122").Value = 1795.9231
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5387.7693
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2937.7693

$ws.Range("H136").Value = 3730.182
$ws.Range("I136").Value = 1237.25
$ws.Range("J136").Value = 5154.7144
$ws.Range("K136").Value = 3711.75
$ws.Range("L136").Value = 15464.1432
$ws.Range("M136").Value = -1161.75
$ws.Range("N136").Value = -20564.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 7512.5713
$ws.Range("I80").Value = 154.5
$ws.Range("J80").Value = 8738.916999999999
$ws.Range("K80").Value = 154.5
$ws.Range("L80").Value = 8738.916999999999
$ws.Range("M80").Value = 843.5
$ws.Range("N80").Value = -10734.917

$ws.Range("H83").Value = 7512.5713
$ws.Range("I83").Value = 154.5
$ws.Range("J83").Value = 8738.916999999999
$ws.Range("K83").Value = 772.5
$ws.Range("L83").Value = 43694.585
$ws.Range("M83").Value = 4219.5
$ws.Range("N83").Value = -53678.585

$ws.Range("H134").Value = 4735.061
$ws.Range("I134").Value = 5149.316
$ws.Range("J134").Value = 3304
$ws.Range("K134").Value = 15447.948
$ws.Range("L134").Value = 9912
$ws.Range("M134").Value = -12912.948
$ws.Range("N134").Value = -14982

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2024.9714
$ws.Range("I31").Value = 1102.7222
$ws.Range("J31").Value = 3001.4707
$ws.Range("K31").Value = 1102.7222
$ws.Range("L31").Value = 3001.4707
$ws.Range("M31").Value = -807.7221999999999
$ws.Range("N31").Value = -3591.4707

$ws.Range("H34").Value = 2024.9714
$ws.Range("I34").Value = 1102.7222
$ws.Range("J34").Value = 3001.4707
$ws.Range("K34").Value = 1102.7222
$ws.Range("L34").Value = 3001.4707
$ws.Range("M34").Value = -900.7221999999999
$ws.Range("N34").Value = -3405.4707

$ws.Range("H134").Value = 1620.225
$ws.Range("I134").Value = 926.86664
$ws.Range("J134").Value = 3700.3
$ws.Range("K134").Value = 2780.59992
$ws.Range("L134").Value = 11100.9
$ws.Range("M134").Value = -245.5999199999997

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 168233.17
$ws.Range("I86").Value = 168233.17
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 504699.51
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -503513.51
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 168233.17
$ws.Range("I89").Value = 168233.17
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 1514098.53
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1508170.53
$ws.Range("N89").ClearContents()

$ws.Range("H98").Value = 459.15384
$ws.Range("I98").Value = 549
$ws.Range("J98").Value = 442.81818
$ws.Range("K98").Value = 1647
$ws.Range("L98").Value = 1328.45454
$ws.Range("M98").Value = -149
$ws.Range("N98").Value = -4324.45454

$ws.Range("H104").Value = 4333.222
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 4333.222
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 12999.666
$ws.Range("N104").Value = -18241.666
$ws.Range("M104").ClearContents()

$ws.Range("H131").Value = 17883416
$ws.Range("I131").Value = 250000510
$ws.Range("J131").Value = 28253.809
$ws.Range("K131").Value = 750001530
$ws.Range("L131").Value = 84761.427
$ws.Range("M131").Value = -749996490
$ws.Range("N131").Value = -94841.427

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4771
$ws.Range("I70").Value = 4799.4
$ws.Range("J70").Value = 4700
$ws.Range("K70").Value = 4799.4
$ws.Range("L70").Value = 4700
$ws.Range("M70").Value = -4529.4

$ws.Range("H73").Value = 4771
$ws.Range("I73").Value = 4799.4
$ws.Range("J73").Value = 4700
$ws.Range("K73").Value = 4799.4
$ws.Range("L73").Value = 4700
$ws.Range("M73").Value = -3863.4

$ws.Range("H126").Value = 1770047.8
$ws.Range("I126").Value = 2418311.8
$ws.Range("J126").Value = 113372.78
$ws.Range("K126").Value = 7254935.399999999
$ws.Range("L126").Value = 340118.34
$ws.Range("M126").Value = -7252465.399999999
$ws.Range("N126").Value = -345058.34

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 42238.168
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 42238.168
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 42238.168
$ws.Range("N127").Value = -52158.168
